# Fix the "empty note" bug: footnote 30 only ever contained a stray
# "a" (with no real content of its own) that belongs at the end of the
# text of footnote 29. Move that trailing "a" onto footnote 29 and
# remove the now pointless footnote 30 (both its reference run in the
# body and its definition) entirely.

$d = $word.ActiveDocument

# Footnotes are 1-based in document order: id 29 -> index 9, id 30 -> index 10.
$fn29 = $d.Footnotes.Item(9)
$fn30 = $d.Footnotes.Item(10)

# Append the stray "a" to the end of footnote 29's text.
$fn29.Range.InsertAfter("a")

# Remove footnote 30 completely - this deletes both its <w:footnote>
# definition and the matching <w:footnoteReference> run in the body.
$fn30.Delete()
